# Insert a new data row before the existing row 121 ("Primera" / Provincia de
# Linares / 2021-03-24 record), shifting all subsequent rows down by one
# (old row 121 -> new row 122, ..., old row 215 -> new row 216).
#
# The new row 121 duplicates the (old) row 121 record, except for the
# Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P), Origen (R) and Precio $/Kg (S) fields,
# which take new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 121, shifting rows 121:215 down to 122:216
$ws.Rows.Item(121).Insert(-4121, $null)

# Duplicate the contents of the row that is now at 122 (the original row 121)
# into the newly created blank row 121, to seed all the unchanged columns.
$src = $ws.Range("A122:T122")
$dst = $ws.Range("A121:T121")
$dst.Value2 = $src.Value()

# Now overwrite the columns that differ for the new record.
$ws.Range("D121").Value = 44830
$ws.Range("M121").Value = 150
$ws.Range("N121").Value = 10000
$ws.Range("O121").Value = 12000
$ws.Range("P121").Value = 11000
$ws.Range("R121").Value = "Provincia de Limarí"
$ws.Range("S121").Value = 5500
